$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4345.4287
$ws.Range("J43").Value = 4862.5
$ws.Range("L43").Value = 4862.5
$ws.Range("N43").Value = -5000.5

$ws.Range("H98").Value = 2371.0908
$ws.Range("I98").Value = 1654.1428
$ws.Range("J98").Value = 3625.75
$ws.Range("K98").Value = 1654.1428
$ws.Range("L98").Value = 3625.75
$ws.Range("M98").Value = -156.1428000000001
$ws.Range("N98").Value = -6621.75

$ws.Range("H122").Value = 2371.0908
$ws.Range("I122").Value = 1654.1428
$ws.Range("J122").Value = 3625.75
$ws.Range("K122").Value = 4962.428400000001
$ws.Range("L122").Value = 10877.25
$ws.Range("M122").Value = -2512.428400000001
$ws.Range("N122").Value = -15777.25

$ws.Range("H132").Value = 63882.562
$ws.Range("I132").Value = 67941.47
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 203824.41
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -201294.41
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3482.2
$ws.Range("I45").Value = 3099.25
$ws.Range("K45").Value = 3099.25
$ws.Range("M45").Value = -2722.25

$ws.Range("H61").Value = 3437.4285
$ws.Range("I61").Value = 3135.2
$ws.Range("K61").Value = 3135.2
$ws.Range("M61").Value = -2923.2

$ws.Range("H74").Value = 2106385.8
$ws.Range("I74").Value = 1091601
$ws.Range("J74").Value = 5556653.5
$ws.Range("K74").Value = 1091601
$ws.Range("L74").Value = 5556653.5
$ws.Range("M74").Value = -1090727
$ws.Range("N74").Value = -5558401.5

$ws.Range("H77").Value = 2106385.8
$ws.Range("I77").Value = 1091601
$ws.Range("J77").Value = 5556653.5
$ws.Range("K77").Value = 5458005
$ws.Range("L77").Value = 27783267.5
$ws.Range("M77").Value = -5453637
$ws.Range("N77").Value = -27792003.5

$ws.Range("H88").Value = 2199.8
$ws.Range("J88").Value = 999.5
$ws.Range("L88").Value = 999.5
$ws.Range("N88").Value = -1811.5

$ws.Range("H91").Value = 2199.8
$ws.Range("J91").Value = 999.5
$ws.Range("L91").Value = 999.5
$ws.Range("N91").Value = -3807.5

$ws.Range("H106").Value = 107000
$ws.Range("J106").Value = 107000
$ws.Range("L106").Value = 107000
$ws.Range("N106").Value = -109524

$ws.Range("H122").Value = 1023
$ws.Range("I122").Value = 1023
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3069
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -619
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2785.8386
$ws.Range("I132").Value = 2656.0715
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 7968.2145
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -5438.2145
$ws.Range("N132").Value = -17051

$ws.Range("H136").Value = 3437.4285
$ws.Range("I136").Value = 3135.2
$ws.Range("K136").Value = 9405.599999999999
$ws.Range("M136").Value = -6855.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3710.5405
$ws.Range("I86").Value = 1412.3334
$ws.Range("K86").Value = 1412.3334
$ws.Range("M86").Value = -289.3334

$ws.Range("H89").Value = 3710.5405
$ws.Range("I89").Value = 1412.3334
$ws.Range("K89").Value = 7061.666999999999
$ws.Range("M89").Value = -1445.666999999999

$ws.Range("H105").Value = 4144.9165
$ws.Range("I105").Value = 2448.25
$ws.Range("J105").Value = 4993.25
$ws.Range("K105").Value = 2448.25
$ws.Range("L105").Value = 4993.25
$ws.Range("M105").Value = -701.25
$ws.Range("N105").Value = -8487.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 96460.3
$ws.Range("I6").Value = 119228.43
$ws.Range("K6").Value = 119228.43
$ws.Range("M6").Value = -119115.43

$ws.Range("H10").Value = 112
$ws.Range("I10").Value = 112
$ws.Range("K10").Value = 112
$ws.Range("M10").Value = 27

$ws.Range("H19").Value = 136.09091
$ws.Range("I19").Value = 129.7
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 129.7
$ws.Range("L19").Value = 200
$ws.Range("M19").Value = 40.30000000000001
$ws.Range("N19").Value = -540

$ws.Range("H24").Value = 136.09091
$ws.Range("I24").Value = 129.7
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 129.7
$ws.Range("L24").Value = 200
$ws.Range("M24").Value = 40.30000000000001
$ws.Range("N24").Value = -540

$ws.Range("H99").Value = 1849.579
$ws.Range("I99").Value = 1348.9
$ws.Range("J99").Value = 2405.889
$ws.Range("K99").Value = 1348.9
$ws.Range("L99").Value = 2405.889
$ws.Range("M99").Value = 149.0999999999999
$ws.Range("N99").Value = -5401.889

$ws.Range("H107").Value = 869.381
$ws.Range("I107").Value = 862.26666
$ws.Range("J107").Value = 887.1667
$ws.Range("K107").Value = 862.26666
$ws.Range("L107").Value = 887.1667
$ws.Range("M107").Value = 1057.73334
$ws.Range("N107").Value = -4727.1667

$ws.Range("H126").Value = 1849.579
$ws.Range("I126").Value = 1348.9
$ws.Range("J126").Value = 2405.889
$ws.Range("K126").Value = 4046.7
$ws.Range("L126").Value = 7217.667
$ws.Range("M126").Value = -1576.7
$ws.Range("N126").Value = -12157.667

$ws.Range("H134").Value = 2501974.5
$ws.Range("I134").Value = 1951.6562
$ws.Range("J134").Value = 12502065
$ws.Range("K134").Value = 5854.9686
$ws.Range("L134").Value = 37506195
$ws.Range("M134").Value = -3319.9686
$ws.Range("N134").Value = -37511265

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1090.1666
$ws.Range("I12").Value = 180
$ws.Range("J12").Value = 1272.2
$ws.Range("K12").Value = 540
$ws.Range("L12").Value = 3816.6
$ws.Range("M12").Value = -367
$ws.Range("N12").Value = -4162.6

$ws.Range("H14").Value = 219.5
$ws.Range("I14").Value = 219.5
$ws.Range("K14").Value = 658.5
$ws.Range("M14").Value = -485.5

$ws.Range("H26").Value = 90.8
$ws.Range("I26").Value = 39.22222
$ws.Range("K26").Value = 117.66666
$ws.Range("M26").Value = 170.33334

$ws.Range("H52").Value = 1966.6666
$ws.Range("J52").Value = 1966.6666
$ws.Range("L52").Value = 5899.9998
$ws.Range("N52").Value = -6431.9998

$ws.Range("H109").Value = 493.94446
$ws.Range("J109").Value = 662.5714
$ws.Range("L109").Value = 1987.7142
$ws.Range("N109").Value = -4067.7142

$ws.Range("H116").Value = 133042.27
$ws.Range("I116").Value = 157718.44
$ws.Range("J116").Value = 21999.5
$ws.Range("K116").Value = 473155.32
$ws.Range("L116").Value = 65998.5
$ws.Range("M116").Value = -469713.32
$ws.Range("N116").Value = -72882.5

$ws.Range("H134").Value = 4230.6
$ws.Range("I134").Value = 4230.6
$ws.Range("K134").Value = 12691.8
$ws.Range("M134").Value = -7621.800000000001

$ws.Range("H139").Value = 2461
$ws.Range("J139").Value = 3249
$ws.Range("L139").Value = 9747
$ws.Range("N139").Value = -20027

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H102").Value = 1717.5714
$ws.Range("I102").Value = 1407.6666
$ws.Range("K102").Value = 1407.6666
$ws.Range("M102").Value = 214.3334

$ws.Range("H122").Value = 3909
$ws.Range("I122").Value = 4570.1
$ws.Range("J122").Value = 2807.1667
$ws.Range("K122").Value = 13710.3
$ws.Range("L122").Value = 8421.500100000001
$ws.Range("M122").Value = -11260.3
$ws.Range("N122").Value = -13321.5001

$ws.Range("H136").Value = 28866
$ws.Range("J136").Value = 28866
$ws.Range("L136").Value = 86598
$ws.Range("N136").Value = -91698

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1575.8948
$ws.Range("I16").Value = 1526
$ws.Range("K16").Value = 1526
$ws.Range("M16").Value = -1356

$ws.Range("H22").Value = 2090
$ws.Range("J22").Value = 2346.75
$ws.Range("L22").Value = 2346.75
$ws.Range("N22").Value = -2936.75

$ws.Range("H27").Value = 2090
$ws.Range("J27").Value = 2346.75
$ws.Range("L27").Value = 2346.75
$ws.Range("N27").Value = -2560.75

$ws.Range("H61").Value = 3148.1875
$ws.Range("I61").Value = 1497.6666
$ws.Range("J61").Value = 5270.2856
$ws.Range("K61").Value = 1497.6666
$ws.Range("L61").Value = 5270.2856
$ws.Range("M61").Value = -1295.6666
$ws.Range("N61").Value = -5674.2856

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 3148.1875
$ws.Range("I113").Value = 1497.6666
$ws.Range("J113").Value = 5270.2856
$ws.Range("K113").Value = 1497.6666
$ws.Range("L113").Value = 5270.2856
$ws.Range("M113").Value = 672.3334
$ws.Range("N113").Value = -9610.285599999999

$ws.Range("H122").Value = 3548.7856
$ws.Range("I122").Value = 3266.1428
$ws.Range("K122").Value = 9798.428400000001
$ws.Range("M122").Value = -7348.428400000001

$ws.Range("H132").Value = 2440.2104
$ws.Range("J132").Value = 3624.5
$ws.Range("L132").Value = 10873.5
$ws.Range("N132").Value = -15933.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1443.2667
$ws.Range("I100").Value = 1465.96
$ws.Range("K100").Value = 2931.92
$ws.Range("M100").Value = -2390.92

$ws.Range("H132").Value = 1287.1852
$ws.Range("I132").Value = 989.6957
$ws.Range("K132").Value = 2969.0871
$ws.Range("M132").Value = -439.0870999999997
